$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.219.79"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.905.37"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'307.10"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.5251"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").Value = "'0.3806"
$ws.Range("E8").Value = "  +1.06%  "

$ws.Range("D9").Value = "'0.07282"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'21.78"
$ws.Range("E10").Value = "  +2.48%  "

$ws.Range("D11").Value = "'0.9023"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "'0.08193"
$ws.Range("E12").Value = "  -3.37%  "

$ws.Range("D13").Value = "'96.29"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").Value = "'5.359"

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.431.45"
$ws.Range("E16").Value = "  -24.88%  "

$ws.Range("D17").Value = "'0.000008668"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "'14.77"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "27.261.15"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "'5.117"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'10.82"
$ws.Range("E22").Value = "  +1.60%  "

$ws.Range("D23").Value = "'6.498"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'150.11"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").Value = "'2.330"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "'1.742"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "'116.61"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").Value = "'4.844"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").Value = "'4.859"
$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("D31").Value = "'0.09240"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").Value = "'0.8322"
$ws.Range("E32").Value = "  +4.17%  "

$ws.Range("D33").Value = "'0.05050"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("D35").Value = "'2.998"
$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.737"
$ws.Range("E36").Value = "  +4.86%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.347"
$ws.Range("E37").Value = "  -2.85%  "

$ws.Range("D38").Value = "'0.5807"
$ws.Range("E38").Value = "  +1.45%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("D40").Value = "'1.077"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").Value = "'9.251"
$ws.Range("E41").Value = "  +1.46%  "

$ws.Range("D42").Value = "'6.615"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "'117.42"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").Value = "'0.1521"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").Value = "'0.4926"
$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("D46").Value = "'10.21"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").Value = "'1.643"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").Value = "'38.87"
$ws.Range("E49").Value = "  +2.77%  "

$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").Value = "'0.06066"
$ws.Range("E51").Value = "  +1.77%  "
